$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "53.352.55"
$ws.Range("E2").Value = "  -9.00%  "

# Row 3
$ws.Range("D3").Value = "2.375.83"
$ws.Range("E3").Value = "  -12.49%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "455.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -9.46%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.32%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.477"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -9.83%  "

# Row 9
$ws.Range("D9").Value = "2.391.98"
$ws.Range("E9").Value = "  -12.37%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0934"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -9.19%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.20"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -14.04%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.310"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -9.92%  "

# Row 13
$ws.Range("E13").Value = "  -4.61%  "

# Row 14
$ws.Range("D14").Value = "2.797.36"
$ws.Range("E14").Value = "  -11.78%  "

# Row 15
$ws.Range("D15").Value = "53.286.38"
$ws.Range("E15").Value = "  -8.97%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -9.41%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000128"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.28%  "

# Row 18
$ws.Range("D18").Value = "2.395.08"
$ws.Range("E18").Value = "  -11.36%  "

# Row 19
$ws.Range("E19").Value = "  -11.65%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "303.48"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -11.03%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -15.33%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.991"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.02%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.32%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.25"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -15.11%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "55.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -11.32%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.378"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -10.13%  "

# Row 28
$ws.Range("D28").Value = "2.512.86"
$ws.Range("E28").Value = "  -10.58%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.150"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -11.60%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.13%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.04%  "

# Row 32
$ws.Range("D32").Value = "0.0₃0702"
$ws.Range("E32").Value = "  -14.98%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "145.09"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.24%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.31%  "

# Row 35
$ws.Range("E35").Value = "  -11.86%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.91"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.62%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -16.55%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.04"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.52%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.776"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -16.36%  "

# Row 40
$ws.Range("E40").Value = "  +0.00%  "

# Row 41
$ws.Range("E41").Value = "  -8.98%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.589"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.48%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0516"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.13%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.94%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.81%  "

# Row 48
$ws.Range("E48").Value = "  -4.86%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0854"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.03%  "

# Row 46 - was Stacks, now Maker
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.929.89"
$ws.Range("E46").Value = "  -10.74%  "

# Row 47 - was Maker, now Stacks
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -11.47%  "

# Row 50 - was RenderToken, now EnergySwap
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -13.21%  "

# Row 51 - was EnergySwap, now RenderToken
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -10.87%  "

